$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update financial data for rows 2-6 (data series realigned to corrected units)
# Row 2
$ws.Range("D2").Value = 7506
$ws.Range("E2").Value = 480
$ws.Range("F2").Value = 480
$ws.Range("G2").Value = 412
$ws.Range("H2").Value = 322
$ws.Range("I2").Value = 315
$ws.Range("J2").Value = 7
$ws.Range("K2").Value = 5301
$ws.Range("L2").Value = 2553
$ws.Range("M2").Value = 2748
$ws.Range("N2").Value = 2679
$ws.Range("O2").Value = 69
$ws.Range("P2").Value = 428
$ws.Range("Q2").Value = 320
$ws.Range("R2").Value = -441
$ws.Range("S2").Value = 102
$ws.Range("T2").Value = 430
$ws.Range("U2").Value = -110
$ws.Range("V2").Value = 1761
$ws.Range("W2").Value = 6.4
$ws.Range("X2").Value = 4.29
$ws.Range("Y2").Value = 12.46
$ws.Range("Z2").Value = 6.33
$ws.Range("AA2").Value = 92.92
$ws.Range("AB2").Value = 526.11
$ws.Range("AC2").Value = 3676
$ws.Range("AD2").Value = 6.61
$ws.Range("AE2").Value = 31481
$ws.Range("AF2").Value = 0.77
$ws.Range("AG2").Value = 800
$ws.Range("AH2").Value = 3.29
$ws.Range("AI2").Value = 21.61
$ws.Range("AJ2").Value = 8570000

# Row 3
$ws.Range("D3").Value = 6813
$ws.Range("E3").Value = 349
$ws.Range("F3").Value = 349
$ws.Range("G3").Value = 259
$ws.Range("H3").Value = 192
$ws.Range("I3").Value = 194
$ws.Range("J3").Value = -2
$ws.Range("K3").Value = 5003
$ws.Range("L3").Value = 2129
$ws.Range("M3").Value = 2874
$ws.Range("N3").Value = 2792
$ws.Range("O3").Value = 82
$ws.Range("P3").Value = 428
$ws.Range("Q3").Value = 701
$ws.Range("R3").Value = -205
$ws.Range("S3").Value = -431
$ws.Range("T3").Value = 235
$ws.Range("U3").Value = 465
$ws.Range("V3").Value = 1412
$ws.Range("W3").Value = 5.12
$ws.Range("X3").Value = 2.82
$ws.Range("Y3").Value = 7.08
$ws.Range("Z3").Value = 3.73
$ws.Range("AA3").Value = 74.09
$ws.Range("AB3").Value = 554.34
$ws.Range("AC3").Value = 2260
$ws.Range("AD3").Value = 8.359999999999999
$ws.Range("AE3").Value = 33120
$ws.Range("AF3").Value = 0.57
$ws.Range("AG3").Value = 700
$ws.Range("AH3").Value = 3.7
$ws.Range("AI3").Value = 30.47
$ws.Range("AJ3").Value = 8570000

# Row 4
$ws.Range("D4").Value = 6688
$ws.Range("E4").Value = 404
$ws.Range("F4").Value = 404
$ws.Range("G4").Value = 314
$ws.Range("H4").Value = 237
$ws.Range("I4").Value = 236
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 5168
$ws.Range("L4").Value = 2110
$ws.Range("M4").Value = 3059
$ws.Range("N4").Value = 2957
$ws.Range("O4").Value = 101
$ws.Range("P4").Value = 428
$ws.Range("Q4").Value = 630
$ws.Range("R4").Value = -227
$ws.Range("S4").Value = -239
$ws.Range("T4").Value = 230
$ws.Range("U4").Value = 400
$ws.Range("V4").Value = 1211
$ws.Range("W4").Value = 6.04
$ws.Range("X4").Value = 3.54
$ws.Range("Y4").Value = 8.210000000000001
$ws.Range("Z4").Value = 4.65
$ws.Range("AA4").Value = 68.98
$ws.Range("AB4").Value = 595.88
$ws.Range("AC4").Value = 2755
$ws.Range("AD4").Value = 8.17
$ws.Range("AE4").Value = 35105
$ws.Range("AF4").Value = 0.64
$ws.Range("AG4").Value = 750
$ws.Range("AH4").Value = 3.33
$ws.Range("AI4").Value = 26.76
$ws.Range("AJ4").Value = 8570000

# Row 5
$ws.Range("D5").Value = 7299
$ws.Range("E5").Value = 351
$ws.Range("F5").Value = 351
$ws.Range("G5").Value = 313
$ws.Range("H5").Value = 237
$ws.Range("I5").Value = 237
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 5273
$ws.Range("L5").Value = 2060
$ws.Range("M5").Value = 3213
$ws.Range("N5").Value = 3116
$ws.Range("O5").Value = 97
$ws.Range("P5").Value = 428
$ws.Range("Q5").Value = 29
$ws.Range("R5").Value = -104
$ws.Range("S5").Value = -73
$ws.Range("T5").Value = 103
$ws.Range("U5").Value = -74
$ws.Range("V5").Value = 1157
$ws.Range("W5").Value = 4.8
$ws.Range("X5").Value = 3.25
$ws.Range("Y5").Value = 7.81
$ws.Range("Z5").Value = 4.54
$ws.Range("AA5").Value = 64.13
$ws.Range("AB5").Value = 637.37
$ws.Range("AC5").Value = 2766
$ws.Range("AD5").Value = 7.27
$ws.Range("AE5").Value = 36995
$ws.Range("AF5").Value = 0.54
$ws.Range("AG5").Value = 750
$ws.Range("AH5").Value = 3.73
$ws.Range("AI5").Value = 26.65
$ws.Range("AJ5").Value = 8570000

# Row 6
$ws.Range("D6").Value = 7236
$ws.Range("E6").Value = 197
$ws.Range("F6").Value = 197
$ws.Range("G6").Value = 153
$ws.Range("H6").Value = 111
$ws.Range("I6").Value = 112
$ws.Range("K6").Value = 5000
$ws.Range("L6").Value = 1749
$ws.Range("M6").Value = 3251
$ws.Range("N6").Value = 3154
$ws.Range("P6").Value = 428
$ws.Range("Q6").Value = 358
$ws.Range("R6").Value = -76
$ws.Range("S6").Value = -310
$ws.Range("T6").Value = 78
$ws.Range("U6").Value = 280
$ws.Range("V6").Value = 920
$ws.Range("W6").Value = 2.72
$ws.Range("X6").Value = 1.53
$ws.Range("Y6").Value = 3.56
$ws.Range("Z6").Value = 2.15
$ws.Range("AA6").Value = 53.81
$ws.Range("AB6").Value = 646.5599999999999
$ws.Range("AC6").Value = 1302
$ws.Range("AD6").Value = 11.44
$ws.Range("AE6").Value = 37475
$ws.Range("AF6").Value = 0.4
$ws.Range("AG6").Value = 900
$ws.Range("AH6").Value = 6.04
$ws.Range("AI6").Value = 67.88
$ws.Range("AJ6").Value = 8570000

# Rows 7-9: clear all trailing financial metric columns (data no longer available / corrected)
$ws.Range("D7:AJ9").ClearContents()

